$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.913.43"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.45"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.69"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5035"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07171"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8945"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.71"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.12"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07489"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.39"
$ws.Range("E14").Value = "  +6.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.234"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008511"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.20"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9990"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.961.37"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.023"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.104.98"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.37"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.94"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.779"
$ws.Range("E26").Value = "  -3.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.87"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.089"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.03"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.684"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09224"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05136"
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7476"
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.973"
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.153"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +6.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.585"
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02003"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5564"
$ws.Range("E40").Value = "  +4.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.067"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.557"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.21"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.593"
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1471"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4693"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9988"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.562"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.67"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.02"
$ws.Range("E51").Value = "  -1.49%  "
